$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New 2D training schedule data (rows 2-6, columns A-J)
# Columns: A=trialTrain, B=x_fixStart, C=y_fixStart, D=x_corrSteps, E=y_corrSteps,
#          F=x_nrSteps, G=y_nrSteps, H=alienID, I=praclen, J=version
$data = @(
    @(1, 0, 8, 2, 4, 2, -4, 45, 5),
    @(2, 1, 7, 2, 2, 1, -5, 56, 5),
    @(3, 1, 9, 6, 8, 5, -1, 12, 5),
    @(4, 0, 6, 3, 3, 3, -3, 34, 5),
    @(5, 3, 9, 7, 7, 4, -2, 23, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$wb.Save()
